$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Change 1: "A diferencia de otros softwares de control o versión control "
# currently spans 3 runs (plain / spellStart+"softwares"+spellEnd / plain)
# with proofErr spell-check markers around "softwares". Collapse it back into
# plain, unmarked text (delete + reinsert so the identical text still counts
# as an edit and the engine merges the run instead of leaving it untouched).
# ---------------------------------------------------------------------------
$find = $d.Content.Find
$ok = $find.Execute("A diferencia de otros softwares de control o versión control ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) { throw "anchor #1 not found" }
$rng = $find.Parent
$rng.Delete()
$rng.InsertAfter("A diferencia de otros softwares de control o versión control ")

# ---------------------------------------------------------------------------
# Change 2: split the run " online donde se puede compartir con otro
# usuarios." into two runs, breaking right after "se pu".
# ---------------------------------------------------------------------------
$find2 = $d.Content.Find
$ok2 = $find2.Execute(" online donde se puede compartir con otro usuarios.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok2) { throw "anchor #2 not found" }
$rng2 = $find2.Parent
$rng2.Delete()
$rng2.InsertAfter(" online donde se pu")
$rng2.Collapse(0)
$rng2.InsertAfter("ede compartir con otro usuarios.")

# ---------------------------------------------------------------------------
# Change 3: the "Soy Ariel y uso git...." paragraph becomes the new
# "Para arrancar a usar git ..." paragraph (keeping the trailing _GoBack
# bookmark), with "git" wrapped in spell-check proofErr markers.
# ---------------------------------------------------------------------------
$find3 = $d.Content.Find
$ok3 = $find3.Execute("Soy Ariel y uso git.aaaaaaaaaaaaaaaaaaaaaaaaaaaa", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok3) { throw "anchor #3 not found" }
$targetPara = $find3.Parent.Paragraphs(1)
$pRng = $targetPara.Range
$xml3 = "<w:p $wns w:rsidR=`"00F64B49`" w:rsidRDefault=`"00F64B49`" w:rsidP=`"004C3A2F`">" + `
        "<w:r><w:t xml:space=`"preserve`">Para arrancar a usar </w:t></w:r>" + `
        "<w:proofErr w:type=`"spellStart`"/>" + `
        "<w:r><w:t>git</w:t></w:r>" + `
        "<w:proofErr w:type=`"spellEnd`"/>" + `
        "<w:r><w:t xml:space=`"preserve`"> lo primero que tenemos que hacer es descargar la aplicación, a continuación le dejamos el link:</w:t></w:r>" + `
        "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" + `
        "<w:bookmarkEnd w:id=`"0`"/>" + `
        "</w:p>"
$pRng.InsertXML($xml3)

# ---------------------------------------------------------------------------
# Change 4: add a new paragraph with the git-scm.com link, followed by a
# blank paragraph, right after the paragraph edited above.
# ---------------------------------------------------------------------------
$anchorPara = $d.Content.Find
$ok4 = $anchorPara.Execute("lo primero que tenemos que hacer es descargar la aplicación, a continuación le dejamos el link:", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok4) { throw "anchor #4 not found" }
$editedPara = $anchorPara.Parent.Paragraphs(1)
$editedPara.Range.InsertParagraphAfter()

$linkPara = $editedPara.Next()
$linkRng = $linkPara.Range
$xmlLink = "<w:p $wns><w:r><w:t>https://git-scm.com/</w:t></w:r></w:p>"
$linkRng.InsertXML($xmlLink)

$linkPara2 = $d.Content.Find
$okLink = $linkPara2.Execute("https://git-scm.com/", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okLink) { throw "link paragraph not found" }
$linkParagraph = $linkPara2.Parent.Paragraphs(1)
$linkParagraph.Range.InsertParagraphAfter()

$blankPara = $linkParagraph.Next()
$blankRng = $blankPara.Range
$xmlBlank = "<w:p $wns/>"
$blankRng.InsertXML($xmlBlank)

Write-Output "done"
